# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new headers, re-using the same formatting as the
# existing header cells (bold, bordered, centered).
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill every data row (2-38) with the team's season record.
$lastRow = 38
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 53
    $ws.Cells.Item($r, 30).Value = 64
    $ws.Cells.Item($r, 31).Value = 0
}
